$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D, shifting existing D:K to F:M
$ws.Range("D:E").Insert()

# Copy number formats from the (now-shifted) old D:K block into the new D:E columns,
# per contiguous data block, so the new columns inherit the correct style (date vs number).
$ws.Range("F7:M35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:M77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:M102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new columns (D = newest quarter, E = next newest quarter) with their values.
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 778500
$ws.Range("E8").Value = 899400
$ws.Range("D9").Value = 455200
$ws.Range("E9").Value = 501300
$ws.Range("D10").Value = 323300
$ws.Range("E10").Value = 398100
$ws.Range("D12").Value = 128600
$ws.Range("E12").Value = 135000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3900
$ws.Range("E14").Value = 31200
$ws.Range("D15").Value = 5500
$ws.Range("E15").Value = 4700
$ws.Range("D17").Value = 730600
$ws.Range("E17").Value = 829400
$ws.Range("D18").Value = 47900
$ws.Range("E18").Value = 70000
$ws.Range("D20").Value = 4300
$ws.Range("E20").Value = -1300
$ws.Range("D21").Value = 82700
$ws.Range("E21").Value = 97300
$ws.Range("D22").Value = 9400
$ws.Range("E22").Value = 14900
$ws.Range("D23").Value = 42800
$ws.Range("E23").Value = 53800
$ws.Range("D24").Value = 9100
$ws.Range("E24").Value = -6100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 33600
$ws.Range("E26").Value = 59900
$ws.Range("D27").Value = 33600
$ws.Range("E27").Value = 59900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 4100
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -4300
$ws.Range("E32").Value = 1300
$ws.Range("D33").Value = 33600
$ws.Range("E33").Value = 64000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 33600
$ws.Range("E35").Value = 64000
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 668800
$ws.Range("E41").Value = 745400
$ws.Range("D42").Value = 119100
$ws.Range("E42").Value = 149000
$ws.Range("D43").Value = 781300
$ws.Range("E43").Value = 812800
$ws.Range("D44").Value = 323100
$ws.Range("E44").Value = 262800
$ws.Range("D45").Value = 197300
$ws.Range("E45").Value = 172600
$ws.Range("D46").Value = 2089700
$ws.Range("E46").Value = 2142600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 67000
$ws.Range("D48").Value = 288700
$ws.Range("E48").Value = 292100
$ws.Range("D49").Value = 437000
$ws.Range("E49").Value = 446200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 802800
$ws.Range("E52").Value = 808600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3618100
$ws.Range("E54").Value = 3756500
$ws.Range("D57").Value = 335500
$ws.Range("E57").Value = 340600
$ws.Range("D58").Value = 10300
$ws.Range("E58").Value = 10500
$ws.Range("D59").Value = 354900
$ws.Range("E59").Value = 611900
$ws.Range("D60").Value = 700800
$ws.Range("E60").Value = 963000
$ws.Range("D61").Value = 684900
$ws.Range("E61").Value = 754700
$ws.Range("D62").Value = 178100
$ws.Range("E62").Value = 109500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1563800
$ws.Range("E66").Value = 1827200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -4864200
$ws.Range("E72").Value = -4947700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2054300
$ws.Range("E76").Value = 1929300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = 33600
$ws.Range("E81").Value = 64000
$ws.Range("D83").Value = 30500
$ws.Range("E83").Value = 28700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -14100
$ws.Range("E89").Value = 68000
$ws.Range("D91").Value = -15300
$ws.Range("E91").Value = -17200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 69400
$ws.Range("E94").Value = -4100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -132600
$ws.Range("E100").Value = -42600
$ws.Range("D101").Value = 700
$ws.Range("E101").Value = -2100
$ws.Range("D102").Value = -76600
$ws.Range("E102").Value = 19200

Write-Output ("Done. UsedRange: " + $ws.UsedRange.Address())
